$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.976.90'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '2.247.88'
$ws.Range("E3").Value = '  +2.15%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'" + '98.37'
$ws.Range("E5").Value = '  +17.63%  '
$ws.Range("D6").Value = "'" + '271.92'
$ws.Range("E6").Value = '  +5.25%  '
$ws.Range("D7").Value = "'" + '0.626'
$ws.Range("E7").Value = '  +0.79%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = "'" + '0.629'
$ws.Range("E9").Value = '  +5.09%  '
$ws.Range("D10").Value = "'" + '48.15'
$ws.Range("E10").Value = '  +8.31%  '
$ws.Range("D11").Value = "'" + '0.0941'
$ws.Range("E11").Value = '  +2.24%  '
$ws.Range("D12").Value = "'" + '8.24'
$ws.Range("E12").Value = '  +14.81%  '
$ws.Range("D13").Value = "'" + '0.104'
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("D14").Value = "'" + '15.36'
$ws.Range("E14").Value = '  +7.03%  '
$ws.Range("D15").Value = '2.581.90'
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").Value = "'" + '0.830'
$ws.Range("E16").Value = '  +6.68%  '
$ws.Range("D17").Value = '2.260.24'
$ws.Range("E17").Value = '  +1.60%  '
$ws.Range("D18").Value = '43.996.54'
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("E19").Value = '  +2.65%  '
$ws.Range("D20").Value = "'" + '6.20'
$ws.Range("E20").Value = '  +4.85%  '
$ws.Range("D21").Value = "'" + '70.94'
$ws.Range("E21").Value = '  +1.82%  '
$ws.Range("D22").Value = "'" + '2.36'
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D23").Value = "'" + '234.89'
$ws.Range("E23").Value = '  +1.99%  '
$ws.Range("D24").Value = "'" + '9.63'
$ws.Range("E24").Value = '  +7.04%  '
$ws.Range("D25").Value = "'" + '0.999'
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = "'" + '11.42'
$ws.Range("E26").Value = '  +7.37%  '
$ws.Range("D27").Value = "'" + '2.50'
$ws.Range("E27").Value = '  +11.81%  '
$ws.Range("E28").Value = '  +0.67%  '
$ws.Range("D29").Value = "'" + '39.84'
$ws.Range("E29").Value = '  +2.02%  '
$ws.Range("D30").Value = "'" + '2.28'
$ws.Range("E30").Value = '  +2.77%  '
$ws.Range("D31").Value = "'" + '173.33'
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("D32").Value = "'" + '0.0911'
$ws.Range("E32").Value = '  +6.43%  '
$ws.Range("D33").Value = "'" + '21.10'
$ws.Range("E33").Value = '  +3.47%  '
$ws.Range("E34").Value = '  +6.41%  '
$ws.Range("E35").Value = '  +1.80%  '
$ws.Range("E36").Value = '  +0.79%  '
$ws.Range("E37").Value = '  -2.73%  '
$ws.Range("E38").Value = '  -2.27%  '
$ws.Range("D39").Value = "'" + '3.56'
$ws.Range("E39").Value = '  +24.71%  '
$ws.Range("D40").Value = "'" + '0.255'
$ws.Range("E40").Value = '  +27.88%  '
$ws.Range("D41").Value = "'" + '12.57'
$ws.Range("E41").Value = '  +0.32%  '
$ws.Range("E42").Value = '  +4.86%  '
$ws.Range("B43").Value = 'THORChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D43").Value = "'" + '5.45'
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").Value = "'" + '62.15'
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("E45").Value = '  +4.69%  '
$ws.Range("D46").Value = "'" + '8.47'
$ws.Range("E46").Value = '  +1.22%  '
$ws.Range("D47").Value = "'" + '100.59'
$ws.Range("E47").Value = '  +1.01%  '
$ws.Range("E48").Value = '  +3.78%  '
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("D50").Value = "'" + '0.432'
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("D51").Value = '2.463.50'
$ws.Range("E51").Value = '  +1.92%  '
